$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the usage date/count header fields for consistency
$ws.Range("K1").Value = "Usage_Date"
$ws.Range("L1").Value = "Usage_Count"

# Update the selection to match the edited cells
$ws.Range("K1:L1").Select()
